$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.064958020203268
$ws.Cells.Item(2, 4).Value = 1.064480623556905
$ws.Cells.Item(2, 5).Value = 1.069849361753674
$ws.Cells.Item(2, 6).Value = 1.080348728488757
$ws.Cells.Item(2, 9).Value = 1.055012512307365
$ws.Cells.Item(2, 10).Value = 1.069915595348227
$ws.Cells.Item(2, 11).Value = 1.067196445904916
$ws.Cells.Item(2, 12).Value = 1.072550770525806
$ws.Cells.Item(2, 13).Value = 1.08302239663545
$ws.Cells.Item(2, 14).Value = 1.071434998214656

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.066365397558997
$ws.Cells.Item(3, 4).Value = 1.065604393985824
$ws.Cells.Item(3, 5).Value = 1.071196844990688
$ws.Cells.Item(3, 6).Value = 1.081946683529981
$ws.Cells.Item(3, 9).Value = 1.055533768733501
$ws.Cells.Item(3, 10).Value = 1.070976219761236
$ws.Cells.Item(3, 11).Value = 1.068134656960493
$ws.Cells.Item(3, 12).Value = 1.073713187149673
$ws.Cells.Item(3, 13).Value = 1.084436701632795
$ws.Cells.Item(3, 14).Value = 1.072497128835986

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.067274332777823
$ws.Cells.Item(4, 4).Value = 1.06633000166604
$ws.Cells.Item(4, 5).Value = 1.072066849455243
$ws.Cells.Item(4, 6).Value = 1.082979134664511
$ws.Cells.Item(4, 9).Value = 1.055868855065068
$ws.Cells.Item(4, 10).Value = 1.071660313729014
$ws.Cells.Item(4, 11).Value = 1.068739609485971
$ws.Cells.Item(4, 12).Value = 1.074462900842979
$ws.Cells.Item(4, 13).Value = 1.085349824793511
$ws.Cells.Item(4, 14).Value = 1.073182194295665

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.067656041891499
$ws.Cells.Item(5, 4).Value = 1.066634682012203
$ws.Cells.Item(5, 5).Value = 1.07243214994586
$ws.Cells.Item(5, 6).Value = 1.083412818872105
$ws.Cells.Item(5, 9).Value = 1.05600920150832
$ws.Cells.Item(5, 10).Value = 1.071947385696139
$ws.Cells.Item(5, 11).Value = 1.06899342592022
$ws.Cells.Item(5, 12).Value = 1.074777501229377
$ws.Cells.Item(5, 13).Value = 1.085733224482051
$ws.Cells.Item(5, 14).Value = 1.073469673937911

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.067720108864071
$ws.Cells.Item(6, 4).Value = 1.066685817922914
$ws.Cells.Item(6, 5).Value = 1.072493459362536
$ws.Cells.Item(6, 6).Value = 1.083485615601862
$ws.Cells.Item(6, 9).Value = 1.056032735666748
$ws.Cells.Item(6, 10).Value = 1.071995555949995
$ws.Cells.Item(6, 11).Value = 1.069036013311411
$ws.Cells.Item(6, 12).Value = 1.07483029024524
$ws.Cells.Item(6, 13).Value = 1.085797571217655
$ws.Cells.Item(6, 14).Value = 1.073517912599053

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.06727943478474
$ws.Cells.Item(7, 4).Value = 1.066334074247971
$ws.Cells.Item(7, 5).Value = 1.072071732375107
$ws.Cells.Item(7, 6).Value = 1.08298493097105
$ws.Cells.Item(7, 9).Value = 1.055870732435568
$ws.Cells.Item(7, 10).Value = 1.071664151639261
$ws.Cells.Item(7, 11).Value = 1.068743002974127
$ws.Cells.Item(7, 12).Value = 1.074467106816493
$ws.Cells.Item(7, 13).Value = 1.08535494966452
$ws.Cells.Item(7, 14).Value = 1.073186037656185

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.065434013067881
$ws.Cells.Item(8, 4).Value = 1.06486073048752
$ws.Cells.Item(8, 5).Value = 1.070305149300568
$ws.Cells.Item(8, 6).Value = 1.080889087447237
$ws.Cells.Item(8, 9).Value = 1.055189131096504
$ws.Cells.Item(8, 10).Value = 1.07027449754229
$ws.Cells.Item(8, 11).Value = 1.06751396333442
$ws.Cells.Item(8, 12).Value = 1.072944125505226
$ws.Cells.Item(8, 13).Value = 1.083500792061199
$ws.Cells.Item(8, 14).Value = 1.071794410091001

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.062168546040504
$ws.Cells.Item(9, 4).Value = 1.062252427949546
$ws.Cells.Item(9, 5).Value = 1.067177271818839
$ws.Cells.Item(9, 6).Value = 1.077183815858336
$ws.Cells.Item(9, 9).Value = 1.053971059045149
$ws.Cells.Item(9, 10).Value = 1.067808632961241
$ws.Cells.Item(9, 11).Value = 1.065331679782367
$ws.Cells.Item(9, 12).Value = 1.070241390601106
$ws.Cells.Item(9, 13).Value = 1.080217651094373
$ws.Cells.Item(9, 14).Value = 1.069325043699408

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.059981957731912
$ws.Cells.Item(10, 4).Value = 1.060505104877066
$ws.Cells.Item(10, 5).Value = 1.06508152414807
$ws.Cells.Item(10, 6).Value = 1.074704934417485
$ws.Cells.Item(10, 9).Value = 1.053147377318359
$ws.Cells.Item(10, 10).Value = 1.066152849521768
$ws.Cells.Item(10, 11).Value = 1.063865375206194
$ws.Cells.Item(10, 12).Value = 1.068426337294965
$ws.Cells.Item(10, 13).Value = 1.078017715025328
$ws.Cells.Item(10, 14).Value = 1.06766690885753

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.059032756698699
$ws.Cells.Item(11, 4).Value = 1.05974641492021
$ws.Cells.Item(11, 5).Value = 1.064171450579313
$ws.Cells.Item(11, 6).Value = 1.073629365575073
$ws.Cells.Item(11, 9).Value = 1.052787909365476
$ws.Cells.Item(11, 10).Value = 1.065432978733627
$ws.Cells.Item(11, 11).Value = 1.063227662697323
$ws.Cells.Item(11, 12).Value = 1.067637166302906
$ws.Cells.Item(11, 13).Value = 1.077062355971542
$ws.Cells.Item(11, 14).Value = 1.0669460157703

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.058679812412419
$ws.Cells.Item(12, 4).Value = 1.059464283556549
$ws.Cells.Item(12, 5).Value = 1.063833009076167
$ws.Cells.Item(12, 6).Value = 1.073229510765094
$ws.Cells.Item(12, 9).Value = 1.052653961068846
$ws.Cells.Item(12, 10).Value = 1.065165142979939
$ws.Cells.Item(12, 11).Value = 1.062990362070303
$ws.Cells.Item(12, 12).Value = 1.067343537883621
$ws.Cells.Item(12, 13).Value = 1.07670706687356
$ws.Cells.Item(12, 14).Value = 1.06667779965912

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.058755536994096
$ws.Cells.Item(13, 4).Value = 1.059524816283523
$ws.Cells.Item(13, 5).Value = 1.063905624156334
$ws.Cells.Item(13, 6).Value = 1.073315296565573
$ws.Cells.Item(13, 9).Value = 1.052682712769922
$ws.Cells.Item(13, 10).Value = 1.065222614851716
$ws.Cells.Item(13, 11).Value = 1.063041283229087
$ws.Cells.Item(13, 12).Value = 1.067406544717366
$ws.Cells.Item(13, 13).Value = 1.076783297034731
$ws.Cells.Item(13, 14).Value = 1.066735353147547

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.059003589789421
$ws.Cells.Item(14, 4).Value = 1.059723100422126
$ws.Cells.Item(14, 5).Value = 1.064143483103669
$ws.Cells.Item(14, 6).Value = 1.073596320462586
$ws.Cells.Item(14, 9).Value = 1.0527768458687
$ws.Cells.Item(14, 10).Value = 1.065410848439419
$ws.Cells.Item(14, 11).Value = 1.063208056075845
$ws.Cells.Item(14, 12).Value = 1.067612905035612
$ws.Cells.Item(14, 13).Value = 1.077032996414846
$ws.Cells.Item(14, 14).Value = 1.066923854048535

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.059156374165467
$ws.Cells.Item(15, 4).Value = 1.059845227207817
$ws.Cells.Item(15, 5).Value = 1.064289982625692
$ws.Cells.Item(15, 6).Value = 1.073769423102103
$ws.Cells.Item(15, 9).Value = 1.052834787823142
$ws.Cells.Item(15, 10).Value = 1.065526766356977
$ws.Cells.Item(15, 11).Value = 1.063310753721035
$ws.Cells.Item(15, 12).Value = 1.067739984574674
$ws.Cells.Item(15, 13).Value = 1.077186787675459
$ws.Cells.Item(15, 14).Value = 1.067039936582833

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.060044901866628
$ws.Cells.Item(16, 4).Value = 1.060555412094982
$ws.Cells.Item(16, 5).Value = 1.065141867176518
$ws.Cells.Item(16, 6).Value = 1.074776269191505
$ws.Cells.Item(16, 9).Value = 1.053171174515823
$ws.Cells.Item(16, 10).Value = 1.066200563147686
$ws.Cells.Item(16, 11).Value = 1.063907638713095
$ws.Cells.Item(16, 12).Value = 1.068478642988244
$ws.Cells.Item(16, 13).Value = 1.078081059824516
$ws.Cells.Item(16, 14).Value = 1.067714690242269

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.060601604505842
$ws.Cells.Item(17, 4).Value = 1.061000328845234
$ws.Cells.Item(17, 5).Value = 1.065675529235069
$ws.Cells.Item(17, 6).Value = 1.075407241530764
$ws.Cells.Item(17, 9).Value = 1.05338142633641
$ws.Cells.Item(17, 10).Value = 1.066622435121914
$ws.Cells.Item(17, 11).Value = 1.064281296787898
$ws.Cells.Item(17, 12).Value = 1.068941110394129
$ws.Cells.Item(17, 13).Value = 1.078641264081827
$ws.Cells.Item(17, 14).Value = 1.06813716132309

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.060926089382773
$ws.Cells.Item(18, 4).Value = 1.061259640284281
$ws.Cells.Item(18, 5).Value = 1.065986554930951
$ws.Cells.Item(18, 6).Value = 1.075775065780761
$ws.Cells.Item(18, 9).Value = 1.053503791973589
$ws.Cells.Item(18, 10).Value = 1.066868226042268
$ws.Cells.Item(18, 11).Value = 1.064498976163734
$ws.Cells.Item(18, 12).Value = 1.069210547666449
$ws.Cells.Item(18, 13).Value = 1.078967754928798
$ws.Cells.Item(18, 14).Value = 1.068383301294744

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.061036691542549
$ws.Cells.Item(19, 4).Value = 1.061348024890753
$ws.Cells.Item(19, 5).Value = 1.066092564458832
$ws.Cells.Item(19, 6).Value = 1.075900448833916
$ws.Cells.Item(19, 9).Value = 1.053545469699059
$ws.Cells.Item(19, 10).Value = 1.066951987139652
$ws.Cells.Item(19, 11).Value = 1.064573153755228
$ws.Cells.Item(19, 12).Value = 1.069302366117606
$ws.Cells.Item(19, 13).Value = 1.079079034905804
$ws.Cells.Item(19, 14).Value = 1.068467181342494

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.060541899473789
$ws.Cells.Item(20, 4).Value = 1.0609526143234
$ws.Cells.Item(20, 5).Value = 1.065618298306388
$ws.Cells.Item(20, 6).Value = 1.0753395660966
$ws.Cells.Item(20, 9).Value = 1.05335889633366
$ws.Cells.Item(20, 10).Value = 1.066577201251166
$ws.Cells.Item(20, 11).Value = 1.064241234687507
$ws.Cells.Item(20, 12).Value = 1.068891524364293
$ws.Cells.Item(20, 13).Value = 1.078581187135367
$ws.Cells.Item(20, 14).Value = 1.068091863215058

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.058930554676831
$ws.Cells.Item(21, 4).Value = 1.059664719575174
$ws.Cells.Item(21, 5).Value = 1.064073450676431
$ws.Cells.Item(21, 6).Value = 1.073513575418504
$ws.Cells.Item(21, 9).Value = 1.052749137806148
$ws.Cells.Item(21, 10).Value = 1.065355430645214
$ws.Cells.Item(21, 11).Value = 1.06315895743358
$ws.Cells.Item(21, 12).Value = 1.067552150801214
$ws.Cells.Item(21, 13).Value = 1.076959478014816
$ws.Cells.Item(21, 14).Value = 1.066868357554706

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.05791530015575
$ws.Cells.Item(22, 4).Value = 1.058853113575681
$ws.Cells.Item(22, 5).Value = 1.063099826900441
$ws.Cells.Item(22, 6).Value = 1.072363526383903
$ws.Cells.Item(22, 9).Value = 1.05236329233646
$ws.Cells.Item(22, 10).Value = 1.064584684329359
$ws.Cells.Item(22, 11).Value = 1.062476020071286
$ws.Cells.Item(22, 12).Value = 1.066707164412964
$ws.Cells.Item(22, 13).Value = 1.07593737435002
$ws.Cells.Item(22, 14).Value = 1.066096516690677

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.058453711512544
$ws.Cells.Item(23, 4).Value = 1.059283539333874
$ws.Cells.Item(23, 5).Value = 1.06361618587271
$ws.Cells.Item(23, 6).Value = 1.072973380021288
$ws.Cells.Item(23, 9).Value = 1.052568071411702
$ws.Cells.Item(23, 10).Value = 1.064993517604655
$ws.Cells.Item(23, 11).Value = 1.062838293985664
$ws.Cells.Item(23, 12).Value = 1.067155382414345
$ws.Cells.Item(23, 13).Value = 1.076479448304742
$ws.Cells.Item(23, 14).Value = 1.066505930556117

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.060568878319769
$ws.Cells.Item(24, 4).Value = 1.060974175083184
$ws.Cells.Item(24, 5).Value = 1.065644159273341
$ws.Cells.Item(24, 6).Value = 1.075370146365139
$ws.Cells.Item(24, 9).Value = 1.053369077508626
$ws.Cells.Item(24, 10).Value = 1.066597641355213
$ws.Cells.Item(24, 11).Value = 1.064259337857751
$ws.Cells.Item(24, 12).Value = 1.068913931121412
$ws.Cells.Item(24, 13).Value = 1.078608334146102
$ws.Cells.Item(24, 14).Value = 1.068112332346396

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.063014406919318
$ws.Cells.Item(25, 4).Value = 1.062928202432302
$ws.Cells.Item(25, 5).Value = 1.067987719455861
$ws.Cells.Item(25, 6).Value = 1.078143210574551
$ws.Cells.Item(25, 9).Value = 1.054287995505564
$ws.Cells.Item(25, 10).Value = 1.068448183909369
$ws.Cells.Item(25, 11).Value = 1.065897847022674
$ws.Cells.Item(25, 12).Value = 1.070942412984806
$ws.Cells.Item(25, 13).Value = 1.0810683546196
$ws.Cells.Item(25, 14).Value = 1.06996550288324
